$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '48.153.40'
$ws.Range('E2').Value = '  +2.11%  '
$ws.Range('D3').Value = '2.522.83'
$ws.Range('E3').Value = '  +1.29%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '323.51'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +0.54%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '109.58'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +0.76%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.535'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +2.05%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.555'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +3.76%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '40.79'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +5.02%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '20.38'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +11.57%  '
$ws.Range('E12').Value = '  +1.68%  '
$ws.Range('E13').Value = '  +1.13%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '7.27'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +1.54%  '
$ws.Range('D15').Value = '2.918.24'
$ws.Range('E15').Value = '  +1.31%  '
$ws.Range('D16').Value = '2.529.35'
$ws.Range('E16').Value = '  +1.70%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.856'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +1.11%  '
$ws.Range('D18').Value = '48.000.36'
$ws.Range('E18').Value = '  +1.93%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '13.21'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +3.79%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '6.62'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +0.13%  '
$ws.Range('D21').Value = '0.0₃0948'
$ws.Range('E21').Value = '  +1.21%  '
$ws.Range('E22').Value = '  -1.31%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '72.09'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +2.08%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '264.42'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +7.27%  '
$ws.Range('E25').Value = '  +0.24%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '26.17'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +1.44%  '
$ws.Range('E27').Value = '  -0.26%  '
$ws.Range('E28').Value = '  +0.79%  '
$ws.Range('E29').Value = '  +2.71%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '2.22'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -3.05%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '36.20'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +2.54%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '49.66'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -0.49%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '19.95'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -0.64%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '5.38'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -0.48%  '
$ws.Range('E35').Value = '  -0.03%  '
$ws.Range('E36').Value = '  +0.92%  '
$ws.Range('E37').Value = '  +1.28%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '4.72'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +0.73%  '
$ws.Range('E39').Value = '  +1.31%  '
$ws.Range('E40').Value = '  +0.74%  '
$ws.Range('B41').Value = 'Monero'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '120.73'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +0.29%  '
$ws.Range('B42').Value = 'EnergySwap'
$ws.Range('C42').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '22.04'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +3.03%  '
$ws.Range('E43').Value = '  -1.10%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.0301'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +1.90%  '
$ws.Range('D45').Value = '2.018.48'
$ws.Range('E45').Value = '  +1.71%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '3.18'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +4.59%  '
$ws.Range('E47').Value = '  +7.07%  '
$ws.Range('E48').Value = '  +0.43%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '9.12'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +0.62%  '
$ws.Range('E50').Value = '  +2.67%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '79.22'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +2.49%  '
